# Update the "想去人数" (column F) counts on the "展览" and "全部类型"
# sheets. Rows are addressed by their worksheet row number (r="N" in the
# underlying XML), which matches the Excel row number one-for-one.

$wb = $excel.ActiveWorkbook

# row -> new value for column F
$updates = @{
    2  = 1076
    5  = 36
    8  = 1915
    9  = 7095
    10 = 492
    11 = 391
    12 = 325
    14 = 388
    16 = 7020
    18 = 1315
    23 = 288
    24 = 121
    29 = 404
    32 = 86
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
